$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.482.02'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '3.391.67'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.48'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.36'
$ws.Range('E6').Value = '  +2.49%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.187'
$ws.Range('E9').Value = '  +5.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.588'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.59'
$ws.Range('E11').Value = '  +7.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000277'
$ws.Range('E12').Value = '  +3.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '701.94'
$ws.Range('E13').Value = '  +6.90%  '
$ws.Range('D14').Value = '3.937.26'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.54'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').Value = '69.358.55'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.386.37'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.120'
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.65'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.33'
$ws.Range('E20').Value = '  +3.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.905'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('E22').Value = '  +3.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.17'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '101.33'
$ws.Range('E24').Value = '  +3.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.95'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.73'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.66'
$ws.Range('E27').Value = '  +4.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '33.50'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.66'
$ws.Range('E29').Value = '  +3.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.06'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.17'
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '556.08'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.107'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '58.35'
$ws.Range('E34').Value = '  +3.83%  '
$ws.Range('D35').Value = '3.732.09'
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.45'
$ws.Range('E36').Value = '  +5.71%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.146'
$ws.Range('E38').Value = '  +11.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.98'
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.22'
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.66'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').Value = '0.0₃0685'
$ws.Range('E42').Value = '  +3.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.340'
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0419'
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('E45').Value = '  -4.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('E46').Value = '  +3.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.130'
$ws.Range('E47').Value = '  +1.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('E49').Value = '  -1.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.29'
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.65'
$ws.Range('E51').Value = '  -1.26%  '
